# Adding a fallback to the entsoe fetching functions
# Shift the timestamp column (A) forward by one day for every data row,
# and refresh the production values (B) for the rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2..97) forward by exactly one day.
$lastRow = 97
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Updated "Actual Production (MW)" values (column B) for rows 25-46.
$newValues = @{
    25 = 1
    26 = 14
    27 = 21
    28 = 40
    29 = 62
    30 = 111
    31 = 212
    32 = 330
    33 = 386
    34 = 469
    35 = 528
    36 = 582
    37 = 633
    38 = 695
    39 = 756
    40 = 799
    41 = 814
    42 = 839
    43 = 839
    44 = 856
    45 = 902
    46 = 961
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
}

Write-Host "Updated $($lastRow - 1) timestamps and $($newValues.Count) production values"
